$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Preserve the two existing cell styles we need to reuse ----
# Style used by the "Date Column" cells (numFmtId 14, short date).
# Style used by the underline-font "missing value" cells (e.g. old B3).
# Stash copies of both on scratch cells far outside the used range before
# clearing the sheet, so PasteSpecial(formats) can reuse the very same
# style entries afterwards instead of Excel minting brand-new ones.
$ws.Cells.Item(2,4).Copy()
$ws.Cells.Item(100,100).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(100,101).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dateStyleCell = $ws.Cells.Item(100,100)
$underlineStyleCell = $ws.Cells.Item(100,101)

# Start clean: clear any old content/formatting in the region we touch.
$ws.Range("A1:J6").Clear()

# ---- Header row ----
# NOTE: values are assigned in a specific order (not strict left-to-right)
# so that brand-new entries land in the shared-strings table in the same
# sequence as the target workbook: Decimal, Nullable Decimal, Float,
# Nullable Float, Nullable Date.
$ws.Cells.Item(1,1).Value = "Int Column"
$ws.Cells.Item(1,2).Value = "Decimal Column"
$ws.Cells.Item(1,5).Value = "Nullable Decimal Column"
$ws.Cells.Item(1,3).Value = "Float Column"
$ws.Cells.Item(1,6).Value = "Nullable Float Column"
$ws.Cells.Item(1,4).Value = "Nullable Int Column"
$ws.Cells.Item(1,7).Value = "String Column"
$ws.Cells.Item(1,8).Value = "Date Column"
$ws.Cells.Item(1,9).Value = "Nullable Date Column"
$ws.Cells.Item(1,10).Value = "Bool column"

# ---- Row 2 ----
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 1.25
$ws.Cells.Item(2,3).Value = 1.25
$ws.Cells.Item(2,4).Value = 1
$ws.Cells.Item(2,5).Value = 1.25
$ws.Cells.Item(2,6).Value = 1.25
$ws.Cells.Item(2,7).Value = "Item 1"
$ws.Cells.Item(2,8).Value = 36526
$ws.Cells.Item(2,9).Value = 36526
$ws.Cells.Item(2,10).Value = 1

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 2.25
$ws.Cells.Item(3,3).Value = 2.25
# D3 stays empty (nullable int with no value)
$ws.Cells.Item(3,7).Value = "Item 2"
$ws.Cells.Item(3,8).Value = 36527
# I3 stays empty (nullable date with no value)
$ws.Cells.Item(3,10).Value = "Y"

# ---- Row 4 ----
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 3.75
$ws.Cells.Item(4,3).Value = 3.75
$ws.Cells.Item(4,4).Value = 3
$ws.Cells.Item(4,5).Value = 3.75
$ws.Cells.Item(4,6).Value = 3.75
$ws.Cells.Item(4,7).Value = "Item 3"
$ws.Cells.Item(4,8).Value = 36528
$ws.Cells.Item(4,9).Value = 36528
$ws.Cells.Item(4,10).Value = 0

# ---- Row 5 ----
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 4.25
$ws.Cells.Item(5,3).Value = 4.25
# D5, E5, F5 stay empty
$ws.Cells.Item(5,7).Value = "Item 4"
$ws.Cells.Item(5,8).Value = 36529
# I5 stays empty (nullable date with no value)
$ws.Cells.Item(5,10).Value = "N"

# ---- Row 6 ----
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 5
$ws.Cells.Item(6,3).Value = 5
$ws.Cells.Item(6,4).Value = 6
$ws.Cells.Item(6,5).Value = 5
$ws.Cells.Item(6,6).Value = 5
$ws.Cells.Item(6,7).Value = "Item 5"
$ws.Cells.Item(6,8).Value = 36530
$ws.Cells.Item(6,9).Value = 36530
$ws.Cells.Item(6,10).Value = "S"

# ---- Number formats for the date columns (H, I) ----
# Paste-special (formats only) so the existing "short date" style is reused
# for every cell, instead of Excel creating a new style/numFmt entry.
$dateStyleCell.Copy()
$ws.Range("H2:H6").PasteSpecial(-4122)
$ws.Range("I2:I6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Empty-but-styled cells ----
# D3: nullable int column, no value, but keeps the underline-font style used
# elsewhere in the sheet for "missing" cells.
$underlineStyleCell.Copy()
$ws.Cells.Item(3,4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clean up the scratch cells used to stash styles.
$dateStyleCell.Clear()
$underlineStyleCell.Clear()

# ---- Column widths (matches Excel's "AutoFit" after adding the new columns) ----
# The host's ColumnWidth setter quantises to 1/6 of a character, while the
# authored workbook's widths are 1/256-quantised "bestFit" values from real
# Excel; these inputs are chosen to land on the nearest representable width.
$ws.Columns.Item(1).ColumnWidth = 9.833333333333334
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 18.0
$ws.Columns.Item(5).ColumnWidth = 22.833333333333332
$ws.Columns.Item(6).ColumnWidth = 18.833333333333332
$ws.Columns.Item(7).ColumnWidth = 12.666666666666666
$ws.Columns.Item(8).ColumnWidth = 11.666666666666666
$ws.Columns.Item(9).ColumnWidth = 19.666666666666668
$ws.Columns.Item(10).ColumnWidth = 11.166666666666666

# ---- Selection, mirroring the saved sheet view ----
$ws.Range("I7").Select()
